$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 188. This pushes the existing rows 188..352
# down to 189..353. Because every row in this block shares identical
# "descriptive" columns (A,B,C,E,F,G,H,I,N,O,Q,R), the shift alone already
# reproduces the desired end state for those columns; only the brand-new
# row 188 needs its values written out explicitly.
$ws.Rows.Item(188).Insert()

$ws.Cells.Item(188, 1).Value = 3
$ws.Cells.Item(188, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(188, 3).Value = "Coquimbo"
$ws.Cells.Item(188, 4).Value = 44790
$ws.Cells.Item(188, 5).Value = 5
$ws.Cells.Item(188, 6).Value = 100112039
$ws.Cells.Item(188, 7).Value = "Ciboulette"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 120
$ws.Cells.Item(188, 11).Value = 1500
$ws.Cells.Item(188, 12).Value = 1500
$ws.Cells.Item(188, 13).Value = 1500
$ws.Cells.Item(188, 14).Value = "$/docena de atados"
$ws.Cells.Item(188, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(188, 16).Value = 500
$ws.Cells.Item(188, 17).Value = 3
$ws.Cells.Item(188, 18).Value = "Hortaliza"
